# Update the "Return_with_prediction" (G) and "return_pct_change" (H) columns
# (and the single affected "mean_return_pct_change" value I2) with the
# refreshed recurrence results for S&P500_returns_compared_annual.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0.04534360234373833
$ws.Range("H2").Value = -5.985862663517158
$ws.Range("I2").Value = -56.8029598418855
$ws.Range("G3").Value = 0.03407287826542454
$ws.Range("H3").Value = -11.16815572100155
$ws.Range("G4").Value = -0.4589643785384184
$ws.Range("H4").Value = -1.368557456725196
$ws.Range("G5").Value = -0.4610416503270645
$ws.Range("H5").Value = 3.715962513778179
$ws.Range("G6").Value = 0.2382153218473858
$ws.Range("H6").Value = 1.961922118118152
$ws.Range("G7").Value = 0.2421637034103828
$ws.Range("H7").Value = 9.787659436313353
$ws.Range("G8").Value = 0.1633661764750221
$ws.Range("H8").Value = -2.062427407849683
$ws.Range("G9").Value = 0.1713597861898632
$ws.Range("H9").Value = -0.3787794477533725
$ws.Range("G10").Value = -0.005012577247821179
$ws.Range("H10").Value = -5.762136270261233
$ws.Range("G11").Value = -0.0208066575765793
$ws.Range("H11").Value = -41.97364585098097
$ws.Range("G12").Value = 0.1336429190736881
$ws.Range("H12").Value = -2.252270821432676
$ws.Range("G13").Value = 0.1301401762167157
$ws.Range("H13").Value = 4.41487549262665
$ws.Range("G14").Value = 0.2581042854651651
$ws.Range("H14").Value = 4.354772370339903
$ws.Range("G15").Value = 0.2546881009330059
$ws.Range("H15").Value = 0.795913716796391
$ws.Range("G16").Value = 0.1377720098139198
$ws.Range("H16").Value = -10.22827634309022
$ws.Range("G17").Value = 0.1396708649494447
$ws.Range("H17").Value = -7.512047511884214
$ws.Range("G18").Value = -0.004264617492738094
$ws.Range("H18").Value = 73.95638529869801
$ws.Range("G19").Value = -0.01435951328948666
$ws.Range("H19").Value = -1605.275937633223
$ws.Range("G20").Value = 0.1351447394418686
$ws.Range("H20").Value = -2.53669286702483
$ws.Range("G21").Value = 0.1459112239959534
$ws.Range("H21").Value = 1.968453918714429
$ws.Range("G22").Value = 0.1701605477493941
$ws.Range("H22").Value = -8.622631096530156
$ws.Range("G23").Value = 0.1722568318382683
$ws.Range("H23").Value = -4.014071185008659
$ws.Range("G24").Value = -0.09542871741270255
$ws.Range("H24").Value = -1.092674403552268
$ws.Range("G25").Value = -0.101974914344438
$ws.Range("H25").Value = -2.375761963133436
$ws.Range("G26").Value = 0.2292643023411338
$ws.Range("H26").Value = -0.3806980946583911
$ws.Range("G27").Value = 0.2412214548313572
$ws.Range("H27").Value = 3.718937677321665
$ws.Range("G28").Value = 0.05996377679399593
$ws.Range("H28").Value = 1.973553958154038
$ws.Range("G29").Value = 0.08177874276486251
$ws.Range("H29").Value = 15.85931466397264
